$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(10, 8).Value = 10000
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 10000
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 10000
$ws.Cells.Item(10, 14).Value = -10586
$ws.Cells.Item(32, 8).Value = 4551.35
$ws.Cells.Item(32, 9).Value = 5479.8
$ws.Cells.Item(32, 10).Value = 4241.8667
$ws.Cells.Item(32, 11).Value = 5479.8
$ws.Cells.Item(32, 12).Value = 4241.8667
$ws.Cells.Item(32, 13).Value = -5153.8
$ws.Cells.Item(32, 14).Value = -4893.8667
$ws.Cells.Item(43, 8).Value = 3329.5
$ws.Cells.Item(43, 9).Value = 3059.5715
$ws.Cells.Item(43, 10).Value = 3959.3333
$ws.Cells.Item(43, 11).Value = 3059.5715
$ws.Cells.Item(43, 12).Value = 3959.3333
$ws.Cells.Item(43, 13).Value = -2990.5715
$ws.Cells.Item(43, 14).Value = -4097.3333
$ws.Cells.Item(51, 8).Value = 5763.467
$ws.Cells.Item(51, 9).Value = 5333
$ws.Cells.Item(51, 10).Value = 5871.0835
$ws.Cells.Item(51, 11).Value = 5333
$ws.Cells.Item(51, 12).Value = 5871.0835
$ws.Cells.Item(51, 13).Value = -4849
$ws.Cells.Item(51, 14).Value = -6839.0835
$ws.Cells.Item(62, 8).Value = 31857.143
$ws.Cells.Item(62, 9).Value = 29666.666
$ws.Cells.Item(62, 10).Value = 45000
$ws.Cells.Item(62, 11).Value = 29666.666
$ws.Cells.Item(62, 12).Value = 45000
$ws.Cells.Item(62, 13).Value = -29042.666
$ws.Cells.Item(62, 14).Value = -46248
$ws.Cells.Item(65, 8).Value = 31857.143
$ws.Cells.Item(65, 9).Value = 29666.666
$ws.Cells.Item(65, 10).Value = 45000
$ws.Cells.Item(65, 11).Value = 148333.33
$ws.Cells.Item(65, 12).Value = 225000
$ws.Cells.Item(65, 13).Value = -145213.33
$ws.Cells.Item(65, 14).Value = -231240
$ws.Cells.Item(103, 8).Value = 204.16667
$ws.Cells.Item(103, 9).Value = 200
$ws.Cells.Item(103, 10).Value = 208.33333
$ws.Cells.Item(103, 11).Value = 600
$ws.Cells.Item(103, 12).Value = 624.99999
$ws.Cells.Item(103, 13).Value = -14
$ws.Cells.Item(103, 14).Value = -1796.99999
$ws.Cells.Item(137, 8).Value = 1580.8182
$ws.Cells.Item(137, 9).Value = 1376.5555
$ws.Cells.Item(137, 10).Value = 2500
$ws.Cells.Item(137, 11).Value = 4129.666499999999
$ws.Cells.Item(137, 12).Value = 7500
$ws.Cells.Item(137, 13).Value = -1579.666499999999
$ws.Cells.Item(137, 14).Value = -12600
$ws.Cells.Item(138, 8).Value = 2426.2341
$ws.Cells.Item(138, 9).Value = 2242.6667
$ws.Cells.Item(138, 10).Value = 2453.0977
$ws.Cells.Item(138, 11).Value = 6728.000100000001
$ws.Cells.Item(138, 12).Value = 7359.293099999999
$ws.Cells.Item(138, 13).Value = -1588.000100000001
$ws.Cells.Item(138, 14).Value = -17639.2931

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5730.4614
$ws.Cells.Item(2, 9).Value = 5388.778
$ws.Cells.Item(2, 10).Value = 6499.25
$ws.Cells.Item(2, 11).Value = 5388.778
$ws.Cells.Item(2, 12).Value = 6499.25
$ws.Cells.Item(2, 13).Value = -5275.778
$ws.Cells.Item(2, 14).Value = -6725.25
$ws.Cells.Item(26, 8).Value = 5500
$ws.Cells.Item(26, 9).Value = 2000
$ws.Cells.Item(26, 10).Value = 12500
$ws.Cells.Item(26, 11).Value = 2000
$ws.Cells.Item(26, 12).Value = 12500
$ws.Cells.Item(26, 13).Value = -1670
$ws.Cells.Item(26, 14).Value = -13160
$ws.Cells.Item(116, 8).Value = 5730.4614
$ws.Cells.Item(116, 9).Value = 5388.778
$ws.Cells.Item(116, 10).Value = 6499.25
$ws.Cells.Item(116, 11).Value = 5388.778
$ws.Cells.Item(116, 12).Value = 6499.25
$ws.Cells.Item(116, 13).Value = -3094.778
$ws.Cells.Item(116, 14).Value = -11087.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5730.4614
$ws.Cells.Item(3, 9).Value = 5388.778
$ws.Cells.Item(3, 10).Value = 6499.25
$ws.Cells.Item(3, 11).Value = 5388.778
$ws.Cells.Item(3, 12).Value = 6499.25
$ws.Cells.Item(3, 13).Value = -5274.778
$ws.Cells.Item(3, 14).Value = -6727.25
$ws.Cells.Item(12, 8).Value = 410
$ws.Cells.Item(12, 9).Value = 221.66667
$ws.Cells.Item(12, 10).Value = 975
$ws.Cells.Item(12, 11).Value = 221.66667
$ws.Cells.Item(12, 12).Value = 975
$ws.Cells.Item(12, 13).Value = -53.66667000000001
$ws.Cells.Item(12, 14).Value = -1311
$ws.Cells.Item(24, 8).Value = 899.6667
$ws.Cells.Item(24, 9).Value = 899.6667
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 899.6667
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = -664.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 4530.5713
$ws.Cells.Item(2, 9).Value = 108
$ws.Cells.Item(2, 10).Value = 6299.6
$ws.Cells.Item(2, 11).Value = 108
$ws.Cells.Item(2, 12).Value = 6299.6
$ws.Cells.Item(2, 13).Value = 5
$ws.Cells.Item(2, 14).Value = -6525.6
$ws.Cells.Item(10, 8).Value = 23020
$ws.Cells.Item(10, 9).Value = 66660
$ws.Cells.Item(10, 10).Value = 1200
$ws.Cells.Item(10, 11).Value = 66660
$ws.Cells.Item(10, 12).Value = 1200
$ws.Cells.Item(10, 13).Value = -66521
$ws.Cells.Item(10, 14).Value = -1478
$ws.Cells.Item(11, 8).Value = 16499.5
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 16499.5
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 16499.5
$ws.Cells.Item(11, 14).Value = -16779.5
$ws.Cells.Item(12, 8).Value = 2855.5
$ws.Cells.Item(12, 9).Value = 1834.8572
$ws.Cells.Item(12, 10).Value = 10000
$ws.Cells.Item(12, 11).Value = 1834.8572
$ws.Cells.Item(12, 12).Value = 10000
$ws.Cells.Item(12, 13).Value = -1664.8572
$ws.Cells.Item(12, 14).Value = -10340
$ws.Cells.Item(14, 8).Value = 5000
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 5000
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 5000
$ws.Cells.Item(14, 14).Value = -5340
$ws.Cells.Item(21, 8).Value = 5000
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 5000
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 5000
$ws.Cells.Item(21, 14).Value = -5470
$ws.Cells.Item(22, 8).Value = 719.5833
$ws.Cells.Item(22, 9).Value = 326.42856
$ws.Cells.Item(22, 10).Value = 1270
$ws.Cells.Item(22, 11).Value = 326.42856
$ws.Cells.Item(22, 12).Value = 1270
$ws.Cells.Item(22, 13).Value = 23.57144
$ws.Cells.Item(22, 14).Value = -1970
$ws.Cells.Item(26, 8).Value = 2000
$ws.Cells.Item(26, 9).Value = 2000
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 2000
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = -1713
$ws.Cells.Item(31, 8).Value = 4143.6206
$ws.Cells.Item(31, 9).Value = 2864.6086
$ws.Cells.Item(31, 10).Value = 9046.5
$ws.Cells.Item(31, 11).Value = 2864.6086
$ws.Cells.Item(31, 12).Value = 9046.5
$ws.Cells.Item(31, 13).Value = -2569.6086
$ws.Cells.Item(31, 14).Value = -9636.5
$ws.Cells.Item(34, 8).Value = 4143.6206
$ws.Cells.Item(34, 9).Value = 2864.6086
$ws.Cells.Item(34, 10).Value = 9046.5
$ws.Cells.Item(34, 11).Value = 2864.6086
$ws.Cells.Item(34, 12).Value = 9046.5
$ws.Cells.Item(34, 13).Value = -2662.6086
$ws.Cells.Item(34, 14).Value = -9450.5
$ws.Cells.Item(58, 8).Value = 6870.4116
$ws.Cells.Item(58, 9).Value = 3483.1667
$ws.Cells.Item(58, 10).Value = 14999.8
$ws.Cells.Item(58, 11).Value = 3483.1667
$ws.Cells.Item(58, 12).Value = 14999.8
$ws.Cells.Item(58, 13).Value = -3280.1667
$ws.Cells.Item(58, 14).Value = -15405.8
$ws.Cells.Item(111, 8).Value = 78639.336
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = 78639.336
$ws.Cells.Item(111, 11).Value = 0
$ws.Cells.Item(111, 12).Value = 78639.336
$ws.Cells.Item(111, 14).Value = -86819.336
$ws.Cells.Item(122, 8).Value = 5711.231
$ws.Cells.Item(122, 9).Value = 5216.4443
$ws.Cells.Item(122, 10).Value = 6824.5
$ws.Cells.Item(122, 11).Value = 15649.3329
$ws.Cells.Item(122, 12).Value = 20473.5
$ws.Cells.Item(122, 13).Value = -13199.3329
$ws.Cells.Item(122, 14).Value = -25373.5
$ws.Cells.Item(132, 8).Value = 4706.923
$ws.Cells.Item(132, 9).Value = 3719.3
$ws.Cells.Item(132, 10).Value = 7999
$ws.Cells.Item(132, 11).Value = 11157.9
$ws.Cells.Item(132, 12).Value = 23997
$ws.Cells.Item(132, 13).Value = -8627.900000000001
$ws.Cells.Item(132, 14).Value = -29057
$ws.Cells.Item(136, 8).Value = 6870.4116
$ws.Cells.Item(136, 9).Value = 3483.1667
$ws.Cells.Item(136, 10).Value = 14999.8
$ws.Cells.Item(136, 11).Value = 10449.5001
$ws.Cells.Item(136, 12).Value = 44999.39999999999
$ws.Cells.Item(136, 13).Value = -7899.500100000001
$ws.Cells.Item(136, 14).Value = -50099.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 300
$ws.Cells.Item(19, 9).Value = 600
$ws.Cells.Item(92, 8).Value = 296
$ws.Cells.Item(92, 9).Value = 370
$ws.Cells.Item(92, 10).Value = 148
$ws.Cells.Item(92, 11).Value = 1110
$ws.Cells.Item(92, 12).Value = 444
$ws.Cells.Item(92, 13).Value = 138
$ws.Cells.Item(92, 14).Value = -2940
$ws.Cells.Item(117, 8).Value = 170000000
$ws.Cells.Item(117, 9).Value = 166666670
$ws.Cells.Item(117, 10).Value = 171666670
$ws.Cells.Item(117, 11).Value = 500000010
$ws.Cells.Item(117, 12).Value = 515000010
$ws.Cells.Item(117, 13).Value = -499996568
$ws.Cells.Item(117, 14).Value = -515006894
$ws.Cells.Item(137, 8).Value = 8121.4546
$ws.Cells.Item(137, 9).Value = 1407.7778
$ws.Cells.Item(137, 10).Value = 38333
$ws.Cells.Item(137, 11).Value = 4223.3334
$ws.Cells.Item(137, 12).Value = 114999
$ws.Cells.Item(137, 13).Value = 876.6665999999996
$ws.Cells.Item(137, 14).Value = -125199
$ws.Cells.Item(139, 8).Value = 3026.7334
$ws.Cells.Item(139, 9).Value = 3026.7334
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 9080.200199999999
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 13).Value = -3940.200199999999
$ws.Cells.Item(140, 8).Value = 1468.6923
$ws.Cells.Item(140, 9).Value = 1144.8182
$ws.Cells.Item(140, 10).Value = 3250
$ws.Cells.Item(140, 11).Value = 3434.4546
$ws.Cells.Item(140, 12).Value = 9750
$ws.Cells.Item(140, 13).Value = 1745.5454
$ws.Cells.Item(140, 14).Value = -20110

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 38000
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 38000
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 38000
$ws.Cells.Item(26, 14).Value = -38560
$ws.Cells.Item(50, 8).Value = 38000
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 38000
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 38000
$ws.Cells.Item(50, 14).Value = -38996
$ws.Cells.Item(122, 8).Value = 3311.4443
$ws.Cells.Item(122, 9).Value = 2560.8
$ws.Cells.Item(122, 10).Value = 4249.75
$ws.Cells.Item(122, 11).Value = 7682.400000000001
$ws.Cells.Item(122, 12).Value = 12749.25
$ws.Cells.Item(122, 13).Value = -5232.400000000001
$ws.Cells.Item(122, 14).Value = -17649.25
$ws.Cells.Item(123, 8).Value = 42780
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 42780
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 42780
$ws.Cells.Item(123, 14).Value = -47680

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4245.9
$ws.Cells.Item(7, 9).Value = 1601.75
$ws.Cells.Item(7, 10).Value = 6008.6665
$ws.Cells.Item(7, 11).Value = 1601.75
$ws.Cells.Item(7, 12).Value = 6008.6665
$ws.Cells.Item(7, 13).Value = -1489.75
$ws.Cells.Item(7, 14).Value = -6232.6665
$ws.Cells.Item(22, 8).Value = 4149.75
$ws.Cells.Item(22, 9).Value = 4571.143
$ws.Cells.Item(22, 10).Value = 1200
$ws.Cells.Item(22, 11).Value = 4571.143
$ws.Cells.Item(22, 12).Value = 1200
$ws.Cells.Item(22, 13).Value = -4276.143
$ws.Cells.Item(22, 14).Value = -1790
$ws.Cells.Item(27, 8).Value = 4149.75
$ws.Cells.Item(27, 9).Value = 4571.143
$ws.Cells.Item(27, 10).Value = 1200
$ws.Cells.Item(27, 11).Value = 4571.143
$ws.Cells.Item(27, 12).Value = 1200
$ws.Cells.Item(27, 13).Value = -4464.143
$ws.Cells.Item(27, 14).Value = -1414
$ws.Cells.Item(46, 8).Value = 12443.36
$ws.Cells.Item(46, 9).Value = 2651.3635
$ws.Cells.Item(46, 10).Value = 20137.072
$ws.Cells.Item(46, 11).Value = 2651.3635
$ws.Cells.Item(46, 12).Value = 20137.072
$ws.Cells.Item(46, 13).Value = -2463.3635
$ws.Cells.Item(46, 14).Value = -20513.072
$ws.Cells.Item(126, 8).Value = 4245.9
$ws.Cells.Item(126, 9).Value = 1601.75
$ws.Cells.Item(126, 10).Value = 6008.6665
$ws.Cells.Item(126, 11).Value = 4805.25
$ws.Cells.Item(126, 12).Value = 18025.9995
$ws.Cells.Item(126, 13).Value = -2335.25
$ws.Cells.Item(126, 14).Value = -22965.9995
$ws.Cells.Item(136, 8).Value = 3044
$ws.Cells.Item(136, 9).Value = 3061
$ws.Cells.Item(136, 10).Value = 2912.25
$ws.Cells.Item(136, 11).Value = 9183
$ws.Cells.Item(136, 12).Value = 8736.75
$ws.Cells.Item(136, 13).Value = -6633
$ws.Cells.Item(136, 14).Value = -13836.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value = 33835
$ws.Cells.Item(17, 9).Value = 33835
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 33835
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).Value = -33663
$ws.Cells.Item(107, 8).Value = 522.63635
$ws.Cells.Item(107, 9).Value = 604.1667
$ws.Cells.Item(107, 10).Value = 424.8
$ws.Cells.Item(107, 11).Value = 1812.5001
$ws.Cells.Item(107, 12).Value = 1274.4
$ws.Cells.Item(107, 13).Value = 107.4999
$ws.Cells.Item(107, 14).Value = -5114.4
$ws.Cells.Item(136, 8).Value = 2914.182
$ws.Cells.Item(136, 9).Value = 2105.6
$ws.Cells.Item(136, 10).Value = 11000
$ws.Cells.Item(136, 11).Value = 6316.799999999999
$ws.Cells.Item(136, 12).Value = 33000
$ws.Cells.Item(136, 13).Value = -3766.799999999999
$ws.Cells.Item(136, 14).Value = -38100

Write-Output "Updated workbook values"